$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "extra study" note for the 13-03-2024 entry (row 11, column D)
$ws.Range("D11").Value = "Some Basics JWT,Bcrypt,AggPipeline"

# Update the active selection to match the saved workbook state
$ws.Range("E14").Select()
